$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values as re-pulled data
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -3
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = -3
